$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.078931194557045
$ws.Cells.Item(2, 4).Value = 1.069706332033083
$ws.Cells.Item(2, 5).Value = 1.091854854545393
$ws.Cells.Item(2, 6).Value = 1.097265365047157
$ws.Cells.Item(2, 9).Value = 1.052424093324181
$ws.Cells.Item(2, 10).Value = 1.083818215486275
$ws.Cells.Item(2, 11).Value = 1.072408122771914
$ws.Cells.Item(2, 12).Value = 1.094498783120843
$ws.Cells.Item(2, 13).Value = 1.099895541633926
$ws.Cells.Item(2, 14).Value = 1.085357361667952
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.081028077441054
$ws.Cells.Item(3, 4).Value = 1.071297156754798
$ws.Cells.Item(3, 5).Value = 1.093925094214011
$ws.Cells.Item(3, 6).Value = 1.099355220336491
$ws.Cells.Item(3, 9).Value = 1.053063776295608
$ws.Cells.Item(3, 10).Value = 1.085570578081424
$ws.Cells.Item(3, 11).Value = 1.073813250644275
$ws.Cells.Item(3, 12).Value = 1.096386434234268
$ws.Cells.Item(3, 13).Value = 1.101803780736628
$ws.Cells.Item(3, 14).Value = 1.087112212818986
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.082378548724132
$ws.Cells.Item(4, 4).Value = 1.07232052982498
$ws.Cells.Item(4, 5).Value = 1.095258710284686
$ws.Cells.Item(4, 6).Value = 1.100701674007758
$ws.Cells.Item(4, 9).Value = 1.053473273224076
$ws.Cells.Item(4, 10).Value = 1.086697923549498
$ws.Cells.Item(4, 11).Value = 1.074715985347663
$ws.Cells.Item(4, 12).Value = 1.09760150843662
$ws.Cells.Item(4, 13).Value = 1.103032321902955
$ws.Cells.Item(4, 14).Value = 1.088241159246932
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.082944799335117
$ws.Cells.Item(5, 4).Value = 1.072749345557698
$ws.Cells.Item(5, 5).Value = 1.095817965645694
$ws.Cells.Item(5, 6).Value = 1.101266360939846
$ws.Cells.Item(5, 9).Value = 1.053644379607494
$ws.Cells.Item(5, 10).Value = 1.087170320157604
$ws.Cells.Item(5, 11).Value = 1.07509396856895
$ws.Cells.Item(5, 12).Value = 1.098110832568895
$ws.Cells.Item(5, 13).Value = 1.103547342252667
$ws.Cells.Item(5, 14).Value = 1.088714226712406
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.083039789034411
$ws.Cells.Item(6, 4).Value = 1.072821263637368
$ws.Cells.Item(6, 5).Value = 1.095911786088227
$ws.Cells.Item(6, 6).Value = 1.101361095396705
$ws.Cells.Item(6, 9).Value = 1.053673048129139
$ws.Cells.Item(6, 10).Value = 1.087249548183546
$ws.Cells.Item(6, 11).Value = 1.07515734480073
$ws.Cells.Item(6, 12).Value = 1.098196263646557
$ws.Cells.Item(6, 13).Value = 1.103633731784149
$ws.Cells.Item(6, 14).Value = 1.088793567251233
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.08238612079287
$ws.Cells.Item(7, 4).Value = 1.072326265189239
$ws.Cells.Item(7, 5).Value = 1.095266188528125
$ws.Cells.Item(7, 6).Value = 1.100709224692074
$ws.Cells.Item(7, 9).Value = 1.053475563650083
$ws.Cells.Item(7, 10).Value = 1.08670424174641
$ws.Cells.Item(7, 11).Value = 1.074721041943409
$ws.Cells.Item(7, 12).Value = 1.097608319876921
$ws.Cells.Item(7, 13).Value = 1.103039209320081
$ws.Cells.Item(7, 14).Value = 1.088247486416408
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.079641181326128
$ws.Cells.Item(8, 4).Value = 1.070245215769969
$ws.Cells.Item(8, 5).Value = 1.09255575659156
$ws.Cells.Item(8, 6).Value = 1.097972866225173
$ws.Cells.Item(8, 9).Value = 1.052641199922584
$ws.Cells.Item(8, 10).Value = 1.084411808330174
$ws.Cells.Item(8, 11).Value = 1.072884348062853
$ws.Cells.Item(8, 12).Value = 1.095138059242187
$ws.Cells.Item(8, 13).Value = 1.100541745479277
$ws.Cells.Item(8, 14).Value = 1.085951797481776
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.074754091844387
$ws.Cells.Item(9, 4).Value = 1.066531061108504
$ws.Cells.Item(9, 5).Value = 1.08773248485356
$ws.Cells.Item(9, 6).Value = 1.09310502195826
$ws.Cells.Item(9, 9).Value = 1.051136532707346
$ws.Cells.Item(9, 10).Value = 1.080320777143827
$ws.Cells.Item(9, 11).Value = 1.069597156873933
$ws.Cells.Item(9, 12).Value = 1.090735052530401
$ws.Cells.Item(9, 13).Value = 1.096091917014929
$ws.Cells.Item(9, 14).Value = 1.081854956561897
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.07146018575394
$ws.Cells.Item(10, 4).Value = 1.064021679323892
$ws.Cells.Item(10, 5).Value = 1.084483239234318
$ws.Cells.Item(10, 6).Value = 1.089826813956528
$ws.Cells.Item(10, 9).Value = 1.050109501139201
$ws.Cells.Item(10, 10).Value = 1.077556977330175
$ws.Cells.Item(10, 11).Value = 1.067370059175899
$ws.Cells.Item(10, 12).Value = 1.087764109778325
$ws.Cells.Item(10, 13).Value = 1.093090511413509
$ws.Cells.Item(10, 14).Value = 1.079087231835496
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.070024907238598
$ws.Cells.Item(11, 4).Value = 1.062926829781533
$ws.Cells.Item(11, 5).Value = 1.083067821320809
$ws.Cells.Item(11, 6).Value = 1.088399034435255
$ws.Cells.Item(11, 9).Value = 1.049658934712691
$ws.Cells.Item(11, 10).Value = 1.0763511625339
$ws.Cells.Item(11, 11).Value = 1.066396899773909
$ws.Cells.Item(11, 12).Value = 1.086468779840261
$ws.Cells.Item(11, 13).Value = 1.091782169625628
$ws.Cells.Item(11, 14).Value = 1.077879704643901
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.069490387847994
$ws.Cells.Item(12, 4).Value = 1.062518879473069
$ws.Cells.Item(12, 5).Value = 1.082540759283317
$ws.Cells.Item(12, 6).Value = 1.08786740776236
$ws.Cells.Item(12, 9).Value = 1.049490679208666
$ws.Cells.Item(12, 10).Value = 1.075901870330168
$ws.Cells.Item(12, 11).Value = 1.066034070863142
$ws.Cells.Item(12, 12).Value = 1.085986262962108
$ws.Cells.Item(12, 13).Value = 1.091294846478407
$ws.Cells.Item(12, 14).Value = 1.077429774393706
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.069605107785405
$ws.Cells.Item(13, 4).Value = 1.06260644438787
$ws.Cells.Item(13, 5).Value = 1.082653875933365
$ws.Cells.Item(13, 6).Value = 1.087981502318384
$ws.Cells.Item(13, 9).Value = 1.049526811340022
$ws.Cells.Item(13, 10).Value = 1.075998308921809
$ws.Cells.Item(13, 11).Value = 1.066111960668942
$ws.Cells.Item(13, 12).Value = 1.086089827244472
$ws.Cells.Item(13, 13).Value = 1.091399440503141
$ws.Cells.Item(13, 14).Value = 1.077526349939209
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.069980752359962
$ws.Cells.Item(14, 4).Value = 1.062893134666407
$ws.Cells.Item(14, 5).Value = 1.083024281222065
$ws.Cells.Item(14, 6).Value = 1.088355116470743
$ws.Cells.Item(14, 9).Value = 1.049645045009721
$ws.Cells.Item(14, 10).Value = 1.076314052648212
$ws.Cells.Item(14, 11).Value = 1.066366936025466
$ws.Cells.Item(14, 12).Value = 1.086428923087634
$ws.Cells.Item(14, 13).Value = 1.091741915036005
$ws.Cells.Item(14, 14).Value = 1.077842542057919
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.070212013531427
$ws.Cells.Item(15, 4).Value = 1.063069604169986
$ws.Cells.Item(15, 5).Value = 1.083252324982793
$ws.Cells.Item(15, 6).Value = 1.088585140900558
$ws.Cells.Item(15, 9).Value = 1.049717773603644
$ws.Cells.Item(15, 10).Value = 1.076508406321436
$ws.Cells.Item(15, 11).Value = 1.066523854398047
$ws.Cells.Item(15, 12).Value = 1.086637668063955
$ws.Cells.Item(15, 13).Value = 1.09195274530219
$ws.Cells.Item(15, 14).Value = 1.078037171735652
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.071555247406285
$ws.Cells.Item(16, 4).Value = 1.064094163921975
$ws.Cells.Item(16, 5).Value = 1.084576993958877
$ws.Cells.Item(16, 6).Value = 1.089921392864552
$ws.Cells.Item(16, 9).Value = 1.05013927908106
$ws.Cells.Item(16, 10).Value = 1.077636809104287
$ws.Cells.Item(16, 11).Value = 1.06743445633499
$ws.Cells.Item(16, 12).Value = 1.087849886009037
$ws.Cells.Item(16, 13).Value = 1.093177154968806
$ws.Cells.Item(16, 14).Value = 1.079167176979885
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.072395385073343
$ws.Cells.Item(17, 4).Value = 1.064734606578346
$ws.Cells.Item(17, 5).Value = 1.085405627434405
$ws.Cells.Item(17, 6).Value = 1.090757340100284
$ws.Cells.Item(17, 9).Value = 1.050402099985674
$ws.Cells.Item(17, 10).Value = 1.078342172485652
$ws.Cells.Item(17, 11).Value = 1.068003272546249
$ws.Cells.Item(17, 12).Value = 1.088607871208409
$ws.Cells.Item(17, 13).Value = 1.093942835611402
$ws.Cells.Item(17, 14).Value = 1.079873542058158
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.07288455807638
$ws.Cells.Item(18, 4).Value = 1.065107370125807
$ws.Cells.Item(18, 5).Value = 1.085888140512836
$ws.Cells.Item(18, 6).Value = 1.091244136465302
$ws.Cells.Item(18, 9).Value = 1.050554835108539
$ws.Cells.Item(18, 10).Value = 1.078752726236446
$ws.Cells.Item(18, 11).Value = 1.068334205324654
$ws.Cells.Item(18, 12).Value = 1.089049135772182
$ws.Cells.Item(18, 13).Value = 1.094388606041602
$ws.Cells.Item(18, 14).Value = 1.080284678842369
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.073051207943152
$ws.Cells.Item(19, 4).Value = 1.065234338833569
$ws.Cells.Item(19, 5).Value = 1.086052527991987
$ws.Cells.Item(19, 6).Value = 1.091409987365466
$ws.Cells.Item(19, 9).Value = 1.050606818659833
$ws.Cells.Item(19, 10).Value = 1.078892567499219
$ws.Cells.Item(19, 11).Value = 1.068446901909343
$ws.Cells.Item(19, 12).Value = 1.089199451629631
$ws.Cells.Item(19, 13).Value = 1.09454046120641
$ws.Cells.Item(19, 14).Value = 1.080424718695775
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.072305336031103
$ws.Cells.Item(20, 4).Value = 1.064665975686601
$ws.Cells.Item(20, 5).Value = 1.085316807470102
$ws.Cells.Item(20, 6).Value = 1.090667733649361
$ws.Cells.Item(20, 9).Value = 1.050373960223866
$ws.Cells.Item(20, 10).Value = 1.078266584204605
$ws.Cells.Item(20, 11).Value = 1.067942331890503
$ws.Cells.Item(20, 12).Value = 1.088526635316904
$ws.Cells.Item(20, 13).Value = 1.093860772287593
$ws.Cells.Item(20, 14).Value = 1.079797846433081
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.069870173187928
$ws.Cells.Item(21, 4).Value = 1.062808746964512
$ws.Cells.Item(21, 5).Value = 1.08291524265972
$ws.Cells.Item(21, 6).Value = 1.088245132231614
$ws.Cells.Item(21, 9).Value = 1.049610252974568
$ws.Cells.Item(21, 10).Value = 1.076221112899697
$ws.Cells.Item(21, 11).Value = 1.066291889766628
$ws.Cells.Item(21, 12).Value = 1.086329106039763
$ws.Cells.Item(21, 13).Value = 1.09164110230491
$ws.Cells.Item(21, 14).Value = 1.077749470324301
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.068331008788447
$ws.Cells.Item(22, 4).Value = 1.061633643827671
$ws.Cells.Item(22, 5).Value = 1.081397667827493
$ws.Cells.Item(22, 6).Value = 1.086714487196137
$ws.Cells.Item(22, 9).Value = 1.049124893454625
$ws.Cells.Item(22, 10).Value = 1.074926931971127
$ws.Cells.Item(22, 11).Value = 1.065246340990317
$ws.Cells.Item(22, 12).Value = 1.084939465379039
$ws.Cells.Item(22, 13).Value = 1.090237696631275
$ws.Cells.Item(22, 14).Value = 1.076453451510368
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.069147729010615
$ws.Cells.Item(23, 4).Value = 1.062257299485673
$ws.Cells.Item(23, 5).Value = 1.08220289832196
$ws.Cells.Item(23, 6).Value = 1.087526631606146
$ws.Cells.Item(23, 9).Value = 1.049382688561986
$ws.Cells.Item(23, 10).Value = 1.075613782892475
$ws.Cells.Item(23, 11).Value = 1.065801360580561
$ws.Cells.Item(23, 12).Value = 1.085676908119022
$ws.Cells.Item(23, 13).Value = 1.090982421685812
$ws.Cells.Item(23, 14).Value = 1.077141277838808
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.072346027988549
$ws.Cells.Item(24, 4).Value = 1.064696989490179
$ws.Cells.Item(24, 5).Value = 1.085356943904888
$ws.Cells.Item(24, 6).Value = 1.09070822541017
$ws.Cells.Item(24, 9).Value = 1.050386677112657
$ws.Cells.Item(24, 10).Value = 1.078300741990523
$ws.Cells.Item(24, 11).Value = 1.067969870969012
$ws.Cells.Item(24, 12).Value = 1.088563344961076
$ws.Cells.Item(24, 13).Value = 1.093897855759532
$ws.Cells.Item(24, 14).Value = 1.079832052726973
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.076023681816433
$ws.Cells.Item(25, 4).Value = 1.067497003861766
$ws.Cells.Item(25, 5).Value = 1.088985210323747
$ws.Cells.Item(25, 6).Value = 1.094369138141114
$ws.Cells.Item(25, 9).Value = 1.051529683020811
$ws.Cells.Item(25, 10).Value = 1.081384693757194
$ws.Cells.Item(25, 11).Value = 1.07045314071996
$ws.Cells.Item(25, 12).Value = 1.091879467701383
$ws.Cells.Item(25, 13).Value = 1.097248304454665
$ws.Cells.Item(25, 14).Value = 1.082920384058887
